$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.767.77"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.874.98"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.02"
$ws.Range("E5").Value = "  -2.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  +0.22%  "
$ws.Range("E7").Value = "  -0.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.97"
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2896"
$ws.Range("E9").Value = "  +2.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06586"
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.878.60"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.88"
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6653"
$ws.Range("E14").Value = "  +0.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "85.13"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.801"
$ws.Range("E16").Value = "  +1.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.756.11"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007795"
$ws.Range("E18").Value = "  +5.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.37%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.76"
$ws.Range("E20").Value = "  +2.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.120.62"
$ws.Range("E21").Value = "  +2.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9989"
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.725"
$ws.Range("E23").Value = "  +2.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.092"
$ws.Range("E24").Value = "  +2.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.553"
$ws.Range("E25").Value = "  +1.96%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "147.08"
$ws.Range("E26").Value = "  +2.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "133.85"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.63"
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.918"
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.167"
$ws.Range("E31").Value = "  -0.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08646"
$ws.Range("E32").Value = "  +1.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.901"
$ws.Range("E33").Value = "  +1.96%  "
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7038"
$ws.Range("E35").Value = "  +3.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.103"
$ws.Range("E36").Value = "  -0.88%  "
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.202"
$ws.Range("E38").Value = "  -4.57%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.684"
$ws.Range("E39").Value = "  -1.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9274"
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01634"
$ws.Range("E41").Value = "  +2.15%  "
$ws.Range("E42").Value = "  -1.80%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.46"
$ws.Range("E44").Value = "  -0.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4154"
$ws.Range("E45").Value = "  +1.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.475"
$ws.Range("E46").Value = "  +2.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1255"
$ws.Range("E47").Value = "  +1.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05696"
$ws.Range("E48").Value = "  +2.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.49"
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.172"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.337"
$ws.Range("E51").Value = "  +1.22%  "
